$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Início da implementação do cadastro de patente
# New entry: date 08/10/2013 with 4:10 (h:mm) worked.
$ws.Range("A5").Value = 41555
$ws.Range("A5").NumberFormat = "m/d/yy"

$ws.Range("B5").Value = 0.17361111111111113
$ws.Range("B5").NumberFormat = "h:mm"

# Move the active selection to the next empty cell, as left by the author.
$ws.Range("C5").Select()
